$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2; 'E'='3'; 'G'='48.746633'; 'H'='146.239899'; 'I'='0.1549390820513319'; 'J'='0.1549390820513319'; 'K'='3'; 'M'='2.248835333333334'; 'N'='6.746506'; 'O'='0.03590294220158827'; 'P'='0.03590294220158827'; 'Q'='109.6231506714327'; 'R'='986.6083560428939'; 'S'='0.005562768907656114'; 'T'='0.005562768907656114' }
    @{ Row=3; 'E'='3'; 'G'='48.746633'; 'H'='146.239899'; 'I'='0.1549390820513319'; 'J'='0.1549390820513319'; 'K'='3'; 'M'='44.29005966666667'; 'N'='132.870179'; 'O'='0.7070964373190639'; 'P'='0.7070964373190639'; 'Q'='2158.991284119103'; 'R'='19430.92155707192'; 'S'='0.1095568729199829'; 'T'='0.1095568729199829' }
    @{ Row=4; 'E'='3'; 'G'='48.746633'; 'H'='146.239899'; 'I'='0.1549390820513319'; 'J'='0.1549390820513319'; 'K'='3'; 'M'='16.09762433333333'; 'N'='48.292873'; 'O'='0.2570006204793478'; 'P'='0.2570006204793479'; 'Q'='784.7049855488696'; 'R'='7062.344869939826'; 'S'='0.03981944022369289'; 'T'='0.0398194402236929' }
    @{ Row=5; 'E'='3'; 'G'='18.65566266666667'; 'H'='55.966988'; 'I'='0.05929622356958762'; 'J'='0.05929622356958761'; 'K'='3'; 'M'='2.248835333333334'; 'N'='6.746506'; 'O'='0.03590294220158827'; 'P'='0.03590294220158827'; 'Q'='41.95351337154757'; 'R'='377.581620343928'; 'S'='0.00212890888759136'; 'T'='0.00212890888759136' }
    @{ Row=6; 'E'='3'; 'G'='18.65566266666667'; 'H'='55.966988'; 'I'='0.05929622356958762'; 'J'='0.05929622356958761'; 'K'='3'; 'M'='44.29005966666667'; 'N'='132.870179'; 'O'='0.7070964373190639'; 'P'='0.7070964373190639'; 'Q'='826.2604126278726'; 'R'='7436.343713650853'; 'S'='0.04192814843253011'; 'T'='0.04192814843253011' }
    @{ Row=7; 'E'='3'; 'G'='18.65566266666667'; 'H'='55.966988'; 'I'='0.05929622356958762'; 'J'='0.05929622356958761'; 'K'='3'; 'M'='16.09762433333333'; 'N'='48.292873'; 'O'='0.2570006204793478'; 'P'='0.2570006204793479'; 'Q'='300.3118492973916'; 'R'='2702.806643676524'; 'S'='0.01523916624946614'; 'T'='0.01523916624946615' }
    @{ Row=8; 'E'='3'; 'G'='97.66137466666667'; 'H'='292.984124'; 'I'='0.310412490288807'; 'J'='0.310412490288807'; 'K'='3'; 'M'='2.248835333333334'; 'N'='6.746506'; 'O'='0.03590294220158827'; 'P'='0.03590294220158827'; 'Q'='219.6243500523049'; 'R'='1976.619150470744'; 'S'='0.01114472169749012'; 'T'='0.01114472169749012' }
    @{ Row=9; 'E'='3'; 'G'='97.66137466666667'; 'H'='292.984124'; 'I'='0.310412490288807'; 'J'='0.310412490288807'; 'K'='3'; 'M'='44.29005966666667'; 'N'='132.870179'; 'O'='0.7070964373190639'; 'P'='0.7070964373190639'; 'Q'='4325.428111115356'; 'R'='38928.8530000382'; 'S'='0.219491565982554'; 'T'='0.219491565982554' }
    @{ Row=10; 'E'='3'; 'G'='97.66137466666667'; 'H'='292.984124'; 'I'='0.310412490288807'; 'J'='0.310412490288807'; 'K'='3'; 'M'='16.09762433333333'; 'N'='48.292873'; 'O'='0.2570006204793478'; 'P'='0.2570006204793479'; 'Q'='1572.116121260917'; 'R'='14149.04509134825'; 'S'='0.07977620260876292'; 'T'='0.07977620260876295' }
    @{ Row=11; 'E'='3'; 'G'='138.7199146666667'; 'H'='416.159744'; 'I'='0.4409152985128724'; 'J'='0.4409152985128724'; 'K'='3'; 'M'='2.248835333333334'; 'N'='6.746506'; 'O'='0.03590294220158827'; 'P'='0.03590294220158827'; 'Q'='311.958245539385'; 'R'='2807.624209854464'; 'S'='0.0158301564783037'; 'T'='0.0158301564783037' }
    @{ Row=12; 'E'='3'; 'G'='138.7199146666667'; 'H'='416.159744'; 'I'='0.4409152985128724'; 'J'='0.4409152985128724'; 'K'='3'; 'M'='44.29005966666667'; 'N'='132.870179'; 'O'='0.7070964373190639'; 'P'='0.7070964373190639'; 'Q'='6143.913297541577'; 'R'='55295.21967787419'; 'S'='0.3117696367379236'; 'T'='0.3117696367379236' }
    @{ Row=13; 'E'='3'; 'G'='138.7199146666667'; 'H'='416.159744'; 'I'='0.4409152985128724'; 'J'='0.4409152985128724'; 'K'='3'; 'M'='16.09762433333333'; 'N'='48.292873'; 'O'='0.2570006204793478'; 'P'='0.2570006204793479'; 'Q'='2233.061073856057'; 'R'='20097.54966470451'; 'S'='0.1133155052966451'; 'T'='0.1133155052966451' }
    @{ Row=14; 'E'='3'; 'G'='10.83447233333333'; 'H'='32.503417'; 'I'='0.03443690557740099'; 'J'='0.03443690557740099'; 'K'='3'; 'M'='2.248835333333334'; 'N'='6.746506'; 'O'='0.03590294220158827'; 'P'='0.03590294220158827'; 'Q'='24.36494420122244'; 'R'='219.284497811002'; 'S'='0.001236386230546981'; 'T'='0.001236386230546981' }
    @{ Row=15; 'E'='3'; 'G'='10.83447233333333'; 'H'='32.503417'; 'I'='0.03443690557740099'; 'J'='0.03443690557740099'; 'K'='3'; 'M'='44.29005966666667'; 'N'='132.870179'; 'O'='0.7070964373190639'; 'P'='0.7070964373190639'; 'Q'='479.8594261001826'; 'R'='4318.734834901643'; 'S'='0.02435021324607324'; 'T'='0.02435021324607324' }
    @{ Row=16; 'E'='3'; 'G'='10.83447233333333'; 'H'='32.503417'; 'I'='0.03443690557740099'; 'J'='0.03443690557740099'; 'K'='3'; 'M'='16.09762433333333'; 'N'='48.292873'; 'O'='0.2570006204793478'; 'P'='0.2570006204793479'; 'Q'='174.4092654718934'; 'R'='1569.683389247041'; 'S'='0.00885030610078077'; 'T'='0.00885030610078077' }
)


foreach ($rowEntry in $data) {
    $r = $rowEntry['Row']
    foreach ($col in @('E','G','H','I','J','K','M','N','O','P','Q','R','S','T')) {
        $cellRef = "$col$r"
        $ws.Range($cellRef).Value = [double]$rowEntry[$col]
    }
}
